$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 123500
$ws.Range("I64").Value = 485000
$ws.Range("J64").Value = 3000
$ws.Range("K64").Value = 485000
$ws.Range("L64").Value = 3000
$ws.Range("M64").Value = -484752
$ws.Range("N64").Value = -3496

$ws.Range("H67").Value = 123500
$ws.Range("I67").Value = 485000
$ws.Range("J67").Value = 3000
$ws.Range("K67").Value = 485000
$ws.Range("L67").Value = 3000
$ws.Range("M67").Value = -484142
$ws.Range("N67").Value = -4716

$ws.Range("H107").Value = 925.65216
$ws.Range("I107").Value = 925.65216
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 925.65216
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 994.34784
$ws.Range("N107").ClearContents()

$ws.Range("H137").Value = 25020622
$ws.Range("I137").Value = 966.4375
$ws.Range("J137").Value = 69500010
$ws.Range("K137").Value = 2899.3125
$ws.Range("L137").Value = 208500030
$ws.Range("M137").Value = -349.3125
$ws.Range("N137").Value = -208505130

$ws.Range("H141").Value = 1618.6863
$ws.Range("I141").Value = 841.381
$ws.Range("K141").Value = 2524.143
$ws.Range("M141").Value = 2655.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 19185.422
$ws.Range("I2").Value = 25743.572
$ws.Range("J2").Value = 822.6
$ws.Range("K2").Value = 25743.572
$ws.Range("L2").Value = 822.6
$ws.Range("M2").Value = -25630.572
$ws.Range("N2").Value = -1048.6

$ws.Range("H32").Value = 10604782
$ws.Range("I32").Value = 2616296
$ws.Range("J32").Value = 45463628
$ws.Range("K32").Value = 2616296
$ws.Range("L32").Value = 45463628
$ws.Range("M32").Value = -2616009
$ws.Range("N32").Value = -45464202

$ws.Range("H74").Value = 56792140
$ws.Range("I74").Value = 47620130
$ws.Range("J74").Value = 88894170
$ws.Range("K74").Value = 47620130
$ws.Range("L74").Value = 88894170
$ws.Range("M74").Value = -47619256
$ws.Range("N74").Value = -88895918

$ws.Range("H77").Value = 56792140
$ws.Range("I77").Value = 47620130
$ws.Range("J77").Value = 88894170
$ws.Range("K77").Value = 238100650
$ws.Range("L77").Value = 444470850
$ws.Range("M77").Value = -238096282
$ws.Range("N77").Value = -444479586

$ws.Range("H110").Value = 4750
$ws.Range("I110").Value = 2000
$ws.Range("J110").Value = 13000
$ws.Range("K110").Value = 2000
$ws.Range("L110").Value = 13000
$ws.Range("M110").Value = 45
$ws.Range("N110").Value = -17090

$ws.Range("H116").Value = 19185.422
$ws.Range("I116").Value = 25743.572
$ws.Range("J116").Value = 822.6
$ws.Range("K116").Value = 25743.572
$ws.Range("L116").Value = 822.6
$ws.Range("M116").Value = -23449.572
$ws.Range("N116").Value = -5410.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 19185.422
$ws.Range("I3").Value = 25743.572
$ws.Range("J3").Value = 822.6
$ws.Range("K3").Value = 25743.572
$ws.Range("L3").Value = 822.6
$ws.Range("M3").Value = -25629.572
$ws.Range("N3").Value = -1050.6

$ws.Range("H107").Value = 1250717.8
$ws.Range("I107").Value = 2000559.8
$ws.Range("J107").Value = 981
$ws.Range("K107").Value = 2000559.8
$ws.Range("L107").Value = 981
$ws.Range("M107").Value = -1998639.8
$ws.Range("N107").Value = -4821

$ws.Range("H134").Value = 9741435
$ws.Range("I134").Value = 14286724
$ws.Range("J134").Value = 1787179.9
$ws.Range("K134").Value = 42860172
$ws.Range("L134").Value = 5361539.699999999
$ws.Range("M134").Value = -42857637
$ws.Range("N134").Value = -5366609.699999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1364.4286
$ws.Range("I16").Value = 1055.8572
$ws.Range("J16").Value = 1673
$ws.Range("K16").Value = 1055.8572
$ws.Range("L16").Value = 1673
$ws.Range("M16").Value = -768.8571999999999
$ws.Range("N16").Value = -2247

$ws.Range("H31").Value = 1456959.2
$ws.Range("I31").Value = 1745.8
$ws.Range("J31").Value = 3478089.2
$ws.Range("K31").Value = 1745.8
$ws.Range("L31").Value = 3478089.2
$ws.Range("M31").Value = -1450.8
$ws.Range("N31").Value = -3478679.2

$ws.Range("H34").Value = 1456959.2
$ws.Range("I34").Value = 1745.8
$ws.Range("J34").Value = 3478089.2
$ws.Range("K34").Value = 1745.8
$ws.Range("L34").Value = 3478089.2
$ws.Range("M34").Value = -1543.8
$ws.Range("N34").Value = -3478493.2

$ws.Range("H113").Value = 1364.4286
$ws.Range("I113").Value = 1055.8572
$ws.Range("J113").Value = 1673
$ws.Range("K113").Value = 1055.8572
$ws.Range("L113").Value = 1673
$ws.Range("M113").Value = 1114.1428
$ws.Range("N113").Value = -6013

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 5206.5
$ws.Range("J105").Value = 4600
$ws.Range("L105").Value = 13800
$ws.Range("N105").Value = -19042

$ws.Range("H114").Value = 1321.7307
$ws.Range("I114").Value = 254
$ws.Range("J114").Value = 1642.05
$ws.Range("K114").Value = 762
$ws.Range("L114").Value = 4926.15
$ws.Range("M114").Value = 2492
$ws.Range("N114").Value = -11434.15

$ws.Range("H121").Value = 2782302.2
$ws.Range("I121").Value = 1657.5
$ws.Range("J121").Value = 3709183.8
$ws.Range("K121").Value = 4972.5
$ws.Range("L121").Value = 11127551.4
$ws.Range("M121").Value = -3662.5
$ws.Range("N121").Value = -11130171.4

$ws.Range("H129").Value = 2487.4517
$ws.Range("I129").Value = 2145.2666
$ws.Range("J129").Value = 2808.25
$ws.Range("K129").Value = 6435.7998
$ws.Range("L129").Value = 8424.75
$ws.Range("M129").Value = -1435.7998
$ws.Range("N129").Value = -18424.75

$ws.Range("H131").Value = 1004.6429
$ws.Range("J131").Value = 1204.2858
$ws.Range("L131").Value = 3612.8574
$ws.Range("N131").Value = -13692.8574

$ws.Range("H139").Value = 40060.46
$ws.Range("I139").Value = 50872
$ws.Range("J139").Value = 4022
$ws.Range("K139").Value = 152616
$ws.Range("L139").Value = 12066
$ws.Range("M139").Value = -147476
$ws.Range("N139").Value = -22346

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 250.33333
$ws.Range("I107").Value = 134
$ws.Range("J107").Value = 366.66666
$ws.Range("K107").Value = 134
$ws.Range("L107").Value = 366.66666
$ws.Range("M107").Value = 1786
$ws.Range("N107").Value = -4206.66666

$ws.Range("H113").Value = 28521.445
$ws.Range("I113").Value = 675
$ws.Range("J113").Value = 50798.6
$ws.Range("K113").Value = 675
$ws.Range("L113").Value = 50798.6
$ws.Range("M113").Value = 1495
$ws.Range("N113").Value = -55138.6

$ws.Range("H132").Value = 9673419
$ws.Range("I132").Value = 9905615
$ws.Range("K132").Value = 29716845
$ws.Range("M132").Value = -29714315

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 979.5185
$ws.Range("I61").Value = 860.4211
$ws.Range("J61").Value = 1262.375
$ws.Range("K61").Value = 860.4211
$ws.Range("L61").Value = 1262.375
$ws.Range("M61").Value = -658.4211
$ws.Range("N61").Value = -1666.375

$ws.Range("H113").Value = 979.5185
$ws.Range("I113").Value = 860.4211
$ws.Range("J113").Value = 1262.375
$ws.Range("K113").Value = 860.4211
$ws.Range("L113").Value = 1262.375
$ws.Range("M113").Value = 1309.5789
$ws.Range("N113").Value = -5602.375

$ws.Range("H136").Value = 4832864
$ws.Range("I136").Value = 6537486
$ws.Range("J136").Value = 3101.6667
$ws.Range("K136").Value = 19612458
$ws.Range("L136").Value = 9305.000100000001
$ws.Range("M136").Value = -19609908
$ws.Range("N136").Value = -14405.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 9117.044
$ws.Range("I107").Value = 11387.444
$ws.Range("J107").Value = 943.6
$ws.Range("K107").Value = 34162.33199999999
$ws.Range("L107").Value = 2830.8
$ws.Range("M107").Value = -32242.33199999999
$ws.Range("N107").Value = -6670.8

$ws.Range("H113").Value = 401
$ws.Range("I113").Value = 413.4
$ws.Range("J113").Value = 370
$ws.Range("K113").Value = 1240.2
$ws.Range("L113").Value = 1110
$ws.Range("M113").Value = 929.8000000000002
$ws.Range("N113").Value = -5450

$ws.Range("H132").Value = 542809.9
$ws.Range("I132").Value = 1620.125
$ws.Range("J132").Value = 1895784.2
$ws.Range("K132").Value = 4860.375
$ws.Range("L132").Value = 5687352.6
$ws.Range("M132").Value = -2330.375
$ws.Range("N132").Value = -5692412.6
